# Add two new senior-editor columns (H = Test1(F), I = test2(F)) and refresh
# the whole working/"off" shift grid for every tracked date.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells H1 / I1, formatted like the rest of row 1 ------------
$ws.Range("H1").Value = "Test1(F)"
$ws.Range("I1").Value = "test2(F)"

$headerSample = $ws.Range("G1")
$newHeaders = $ws.Range("H1:I1")
$newHeaders.Font.Bold = $headerSample.Font.Bold
$newHeaders.HorizontalAlignment = $headerSample.HorizontalAlignment
$newHeaders.VerticalAlignment = $headerSample.VerticalAlignment
$newHeaders.Borders.LineStyle = $headerSample.Borders.LineStyle

# --- Full shift grid (columns B..I, rows 2..31) -----------------------------
# Columns B-G carry this sprint's roster updates; columns H/I are the new
# Test1(F) / test2(F) schedules for every date row.
$data = @(
    @("15-24", "15-24", "7-16", "off", "off", "10-19", "off", "off"),
    @("10-19", "off", "off", "off", "off", "7-16", "15-24", "15-24"),
    @("7-16", "15-24", "15-24", "off", "off", "off", "10-19", "off"),
    @("15-24", "10-19", "off", "off", "off", "15-24", "7-16", "off"),
    @("off", "7-16", "off", "10-19", "15-24", "15-24", "off", "off"),
    @("off", "off", "off", "7-16", "10-19", "off", "15-24", "off"),
    @("off", "off", "10-19", "7-16", "15-24", "off", "off", "off"),
    @("off", "7-16", "15-24", "15-24", "off", "10-19", "off", "off"),
    @("15-24", "10-19", "off", "15-24", "off", "off", "7-16", "off"),
    @("15-24", "7-16", "off", "off", "15-24", "off", "10-19", "off"),
    @("off", "15-24", "off", "10-19", "off", "7-16", "15-24", "off"),
    @("10-19", "15-24", "off", "off", "7-16", "off", "15-24", "off"),
    @("15-24", "10-19", "off", "off", "off", "7-16", "off", "off"),
    @("10-19", "off", "15-24", "off", "7-16", "off", "off", "off"),
    @("10-19", "15-24", "15-24", "off", "off", "7-16", "off", "off"),
    @("7-16", "15-24", "15-24", "off", "off", "10-19", "off", "off"),
    @("15-24", "15-24", "7-16", "off", "off", "10-19", "off", "off"),
    @("7-16", "off", "off", "15-24", "10-19", "15-24", "off", "off"),
    @("7-16", "off", "10-19", "off", "15-24", "15-24", "off", "off"),
    @("7-16", "15-24", "off", "off", "10-19", "off", "off", "off"),
    @("off", "15-24", "off", "7-16", "10-19", "off", "off", "off"),
    @("off", "15-24", "10-19", "7-16", "15-24", "off", "off", "off"),
    @("10-19", "15-24", "off", "7-16", "off", "15-24", "off", "off"),
    @("15-24", "15-24", "10-19", "7-16", "off", "off", "off", "off"),
    @("off", "15-24", "10-19", "7-16", "off", "15-24", "off", "off"),
    @("off", "15-24", "15-24", "7-16", "off", "10-19", "off", "off"),
    @("10-19", "off", "15-24", "7-16", "off", "off", "off", "off"),
    @("10-19", "off", "off", "15-24", "off", "7-16", "off", "off"),
    @("off", "15-24", "10-19", "15-24", "off", "7-16", "off", "off"),
    @("off", "15-24", "10-19", "15-24", "off", "7-16", "off", "off")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $ws.Cells.Item($row, 2 + $j).Value = $rowValues[$j]
    }
}
